$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A54").Value = "Bridge Domain Name"
$ws.Range("B54").Value = "Bridge Domain Name"
$ws.Range("C54").Value = "브릿지도메인 이름"
$ws.Range("D54").Value = "브릿지도메인 이름"

$ws.Range("A55").Value = "Application Profile Details"
$ws.Range("B55").Value = "Application Profile Details"
$ws.Range("C55").Value = "프로파일 상세정보"
$ws.Range("D55").Value = "프로파일 상세정보"

$ws.Range("A56").Value = "Bridge Domain Details"
$ws.Range("B56").Value = "Bridge Domain Details"
$ws.Range("C56").Value = "브릿지도메인 상세정보"
$ws.Range("D56").Value = "브릿지도메인 상세정보"

$ws.Range("A57").Value = "Subnet IP"
$ws.Range("B57").Value = "Subnet IP"
$ws.Range("C57").Value = "서브넷 주소"
$ws.Range("D57").Value = "서브넷 주소"

$ws.Range("A58").Value = "Context Name"
$ws.Range("B58").Value = "Context Name"
$ws.Range("C58").Value = "VRF 이름"
$ws.Range("D58").Value = "VRF 이름"

$ws.Range("A59").Value = "Context Details"
$ws.Range("B59").Value = "Context Details"
$ws.Range("C59").Value = "VRF 상세정보"
$ws.Range("D59").Value = "VRF 상세정보"

$ws.Range("A60").Value = "Contract Details"
$ws.Range("B60").Value = "Contract Details"
$ws.Range("C60").Value = "컨트랙 상세정보"
$ws.Range("D60").Value = "컨트랙 상세정보"

$ws.Range("A61").Value = "Contract Name"
$ws.Range("B61").Value = "Contract Name"
$ws.Range("C61").Value = "컨트랙 이름"
$ws.Range("D61").Value = "컨트랙 이름"

$ws.Range("A62").Value = "Filter Name"
$ws.Range("B62").Value = "Filter Name"
$ws.Range("C62").Value = "필터 이름"
$ws.Range("D62").Value = "필터 이름"

$ws.Range("A63").Value = "Filter Details"
$ws.Range("B63").Value = "Filter Details"
$ws.Range("C63").Value = "필터 상세정보"
$ws.Range("D63").Value = "필터 상세정보"

$ws.Range("A64").Value = "L3 External Name"
$ws.Range("B64").Value = "L3 External Name"
$ws.Range("C64").Value = "외부네트워크 이름"
$ws.Range("D64").Value = "외부네트워크 이름"

$ws.Range("A65").Value = "L3 External Details"
$ws.Range("B65").Value = "L3 External Details"
$ws.Range("C65").Value = "외부네트워크 상세정보"
$ws.Range("D65").Value = "외부네트워크 상세정보"

$ws.Range("A66").Value = "Filter Entry Name"
$ws.Range("B66").Value = "Filter Entry Name"
$ws.Range("C66").Value = "필터 엔트리 이름"
$ws.Range("D66").Value = "필터 엔트리 이름"

$ws.Range("A67").Value = "Filter Entry Details"
$ws.Range("B67").Value = "Filter Entry Details"
$ws.Range("C67").Value = "필터 엔트리 상세정보"
$ws.Range("D67").Value = "필터 엔트리 상세정보"

$ws.Range("A68").Value = "Subject Name"
$ws.Range("B68").Value = "Subject Name"
$ws.Range("C68").Value = "서브젝트 이름"
$ws.Range("D68").Value = "서브젝트 이름"

$ws.Range("A69").Value = "Refresh"
$ws.Range("B69").Value = "Refresh"
$ws.Range("C69").Value = "정보갱신"
$ws.Range("D69").Value = "정보갱신"

$ws.Range("A70").Value = "Health"
$ws.Range("B70").Value = "Health"
$ws.Range("C70").Value = "상태점수"
$ws.Range("D70").Value = "상태점수"

$ws.Range("A71").Value = "Tenant Name"
$ws.Range("B71").Value = "Tenant Name"
$ws.Range("C71").Value = "테넌트 이름"
$ws.Range("D71").Value = "테넌트 이름"

$ws.Range("A72").Value = "Device Name"
$ws.Range("B72").Value = "Device Name"
$ws.Range("C72").Value = "장치 이름"
$ws.Range("D72").Value = "장치 이름"

$ws.Range("A73").Value = "EPG Name"
$ws.Range("B73").Value = "EPG Name"
$ws.Range("C73").Value = "엔드포인트그룹 이름"
$ws.Range("D73").Value = "엔드포인트그룹 이름"

$ws.Range("A74").Value = "Provided Contracts"
$ws.Range("B74").Value = "Provided Contracts"
$ws.Range("C74").Value = "제공계약"
$ws.Range("D74").Value = "제공계약"

$ws.Range("A75").Value = "Consumed Contracts"
$ws.Range("B75").Value = "Consumed Contracts"
$ws.Range("C75").Value = "이용계약"
$ws.Range("D75").Value = "이용계약"

$ws.Range("A76").Value = "Encap"
$ws.Range("B76").Value = "Encap"
$ws.Range("C76").Value = "망분리"
$ws.Range("D76").Value = "망분리"

$ws.Range("A77").Value = "Binding Path"
$ws.Range("B77").Value = "Binding Path"
$ws.Range("C77").Value = "할당경로"
$ws.Range("D77").Value = "할당경로"

$ws.Range("A78").Value = "MAC"
$ws.Range("B78").Value = "MAC"
$ws.Range("C78").Value = "MAC"
$ws.Range("D78").Value = "MAC"

$ws.Range("A79").Value = "IP"
$ws.Range("B79").Value = "IP"
$ws.Range("C79").Value = "IP"
$ws.Range("D79").Value = "IP"

$ws.Range("A80").Value = "Endpoint"
$ws.Range("B80").Value = "Endpoint"
$ws.Range("C80").Value = "엔드포인트"
$ws.Range("D80").Value = "엔드포인트"

$ws.Range("A81").Value = "Relations"
$ws.Range("B81").Value = "Relations"
$ws.Range("C81").Value = "연관정보"
$ws.Range("D81").Value = "연관정보"

$ws.Range("A82").Value = "Bridge Domain Relations"
$ws.Range("B82").Value = "Bridge Domain Relations"
$ws.Range("C82").Value = "브릿지도메인 연결"
$ws.Range("D82").Value = "브릿지도메인 연결"

$ws.Range("A83").Value = "Path Name"
$ws.Range("B83").Value = "Path Name"
$ws.Range("C83").Value = "경로 이름"
$ws.Range("D83").Value = "경로 이름"

$ws.Range("A84").Value = "Path Attachments"
$ws.Range("B84").Value = "Path Attachments"
$ws.Range("C84").Value = "경로 설정"
$ws.Range("D84").Value = "경로 설정"

$ws.Range("A85").Value = "Current Health"
$ws.Range("B85").Value = "Current Health"
$ws.Range("C85").Value = "현재상태수치"
$ws.Range("D85").Value = "현재상태수치"

$ws.Range("A86").Value = "Total Health"
$ws.Range("B86").Value = "Total Health"
$ws.Range("C86").Value = "전체 상태정보"
$ws.Range("D86").Value = "전체 상태정보"

$ws.Range("A87").Value = "Node Health"
$ws.Range("B87").Value = "Node Health"
$ws.Range("C87").Value = "노드 상태정보"
$ws.Range("D87").Value = "노드 상태정보"

$ws.Range("A88").Value = "Endpoint Group Health"
$ws.Range("B88").Value = "Endpoint Group Health"
$ws.Range("C88").Value = "엔드포인트그룹 상태정보"
$ws.Range("D88").Value = "엔드포인트그룹 상태정보"

# Update sheet view: scroll position and active selection to match the
# post-edit cursor location (best effort; plain scroll w/o freeze panes
# may not round-trip through every host, selection always will).
try {
    $win = $wb.Windows.Item(1)
    $win.ScrollRow = 34
    $win.ScrollColumn = 1
} catch {
}
$ws.Range("D44").Select() | Out-Null
